# "update plots for each sample"
#
# The CYP2D6_10B marker (peak_table row 3 / allele_table row 4 / marker_table
# row 3) was re-analysed: its "wild-type" peak -- previously reported as
# undetected -- is now detected, so the related height/peak/size/status
# fields, the resulting genotype call and phenotype, and the sample-level
# genotype summary all need to be refreshed together.

$wb = $excel.ActiveWorkbook

# --- peak_table: lower the w_height threshold used for CYP2D6_10B ---------
$peakTable = $wb.Worksheets.Item("peak_table")
$peakTable.Range("N3").Value = 800

# --- allele_table: CYP2D6_10B / wildtype ("C") row now detects a peak -----
$alleleTable = $wb.Worksheets.Item("allele_table")
$alleleTable.Range("K4").Value = 800
$alleleTable.Range("M4").Value = $true
$alleleTable.Range("N4").Value = 46
$alleleTable.Range("O4").Value = 32.95
$alleleTable.Range("P4").Value = 933
$alleleTable.Range("Q4").Value = "ok"
$alleleTable.Range("R4").Value = ""

# --- marker_table: genotype/phenotype call for CYP2D6_10B updates ---------
$markerTable = $wb.Worksheets.Item("marker_table")
$markerTable.Range("G3").Value = "CT"
$markerTable.Range("H3").Value = "heterozygous"

# --- genotype_result: overall sample genotype reflects the new call ------
$genotypeResult = $wb.Worksheets.Item("genotype_result")
$genotypeResult.Range("B2").Value = "*1/*10B"
